$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 113-114 (existing rows 113..142 shift down to 115..144).
$ws.Range("A113:A114").EntireRow.Insert()

# New row 113: 2020-06-01 (serial 43983)
$ws.Cells.Item(113, 1).Value = 43983
$ws.Cells.Item(113, 2).Value = 47669
$ws.Cells.Item(113, 3).Value = 31779
$ws.Cells.Item(113, 4).Value = 15890
$ws.Cells.Item(113, 5).Value = 6355.8
$ws.Cells.Item(113, 6).Value = 7945

# New row 114: 2020-07-01 (serial 44013)
$ws.Cells.Item(114, 1).Value = 44013
$ws.Cells.Item(114, 2).Value = 179843
$ws.Cells.Item(114, 3).Value = 110279
$ws.Cells.Item(114, 4).Value = 69564
$ws.Cells.Item(114, 5).Value = 22055.8
$ws.Cells.Item(114, 6).Value = 34782

# Column A on these rows uses the same bordered/centered date style as the rest
# of column A; forcing the border (present on every other A-column cell) makes
# the engine re-use the existing shared style instead of minting a near-duplicate.
$ws.Cells.Item(113, 1).Borders.LineStyle = 1
$ws.Cells.Item(114, 1).Borders.LineStyle = 1
